$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh: update price (D) and volume-change (E) columns,
# plus the two row-pair reorderings (Cosmos/Maker, ThetaToken/USDe).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.238.58'
$ws.Range("E2").Value = '  +1.80%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.325.29'
$ws.Range("E3").Value = '  +6.31%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '600.80'
$ws.Range("E5").Value = '  +1.10%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.12'
$ws.Range("E6").Value = '  +5.76%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.315.82'
$ws.Range("E8").Value = '  +6.36%  '

$ws.Range("E9").Value = '  +0.74%  '

$ws.Range("E11").Value = '  +2.74%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.475'
$ws.Range("E12").Value = '  +3.40%  '

$ws.Range("E13").Value = '  +0.19%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.10'
$ws.Range("E14").Value = '  +2.94%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.861.86'
$ws.Range("E15").Value = '  +6.10%  '

$ws.Range("E16").Value = '  +1.35%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.314.50'
$ws.Range("E17").Value = '  +6.09%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.255.06'
$ws.Range("E18").Value = '  +1.77%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.93'
$ws.Range("E19").Value = '  +3.16%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '486.21'
$ws.Range("E20").Value = '  +1.97%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.39'
$ws.Range("E21").Value = '  +1.52%  '

$ws.Range("E22").Value = '  +6.91%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.10'
$ws.Range("E23").Value = '  +5.80%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.63'
$ws.Range("E24").Value = '  +4.26%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.09'
$ws.Range("E25").Value = '  -2.98%  '

$ws.Range("E27").Value = '  +2.84%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.34'
$ws.Range("E28").Value = '  +4.44%  '

$ws.Range("E29").Value = '  +0.02%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.24'
$ws.Range("E30").Value = '  +0.91%  '

$ws.Range("E31").Value = '  +5.09%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.62'
$ws.Range("E32").Value = '  +4.61%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.107'
$ws.Range("E33").Value = '  -1.37%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.59'
$ws.Range("E34").Value = '  +1.91%  '

$ws.Range("E35").Value = '  +2.69%  '

$ws.Range("E36").Value = '  +3.26%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '53.37'
$ws.Range("E37").Value = '  +2.51%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0740'
$ws.Range("E38").Value = '  +3.87%  '

$ws.Range("E39").Value = '  +3.05%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '430.90'
$ws.Range("E40").Value = '  +2.50%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.80'
$ws.Range("E41").Value = '  +3.83%  '

$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.025.31'
$ws.Range("E42").Value = '  +5.53%  '

$ws.Range("B43").Value = 'Cosmos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.50'
$ws.Range("E43").Value = '  +2.82%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.112'
$ws.Range("E44").Value = '  -4.85%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.272'
$ws.Range("E45").Value = '  +5.09%  '

$ws.Range("E46").Value = '  +6.50%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '26.46'
$ws.Range("E47").Value = '  +3.73%  '

$ws.Range("B48").Value = 'USDe'
$ws.Range("C48").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.999'
$ws.Range("E48").Value = '  +0.10%  '

$ws.Range("B49").Value = 'ThetaToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.35'
$ws.Range("E49").Value = '  +2.45%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.115'

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '35.32'
$ws.Range("E51").Value = '  +15.13%  '
